# MasterExecutor_Sanity.xlsx — "Enabling all testcases LogixalQA"
#
# The RunMode column (E2:E31) on the MasterExecutor sheet currently reads
# "No" for every testcase row. Flip the whole column to "Yes" so every
# testcase is enabled, then leave the view parked on the top data row
# (E2) instead of the bottom one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")
$ws.Activate()

# Flip the entire RunMode column in one shot so every row that used to
# share the "No" string now shares "Yes" (header in E1 is untouched).
$ws.Range("E2:E31").Value = "Yes"

# Move the selection/scroll position up to the top of the data instead of
# leaving it parked at the last row.
$ws.Range("E2").Select()
